$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values for rows 34-49: Date (col B), Visitor (col C), Home (col D)
$rows = @(
    @{ Row = 34; Date = "01/03/2026"; Visitor = "Carolina Panthers"; Home = "Tampa Bay Buccaneers" },
    @{ Row = 35; Date = "01/03/2026"; Visitor = "Seattle Seahawks"; Home = "San Francisco 49ers" },
    @{ Row = 36; Date = "01/04/2026"; Visitor = "Indianapolis Colts"; Home = "Houston Texans" },
    @{ Row = 37; Date = "01/04/2026"; Visitor = "Green Bay Packers"; Home = "Minnesota Vikings" },
    @{ Row = 38; Date = "01/04/2026"; Visitor = "Cleveland Browns"; Home = "Cincinnati Bengals" },
    @{ Row = 39; Date = "01/04/2026"; Visitor = "New Orleans Saints"; Home = "Atlanta Falcons" },
    @{ Row = 40; Date = "01/04/2026"; Visitor = "Dallas Cowboys"; Home = "New York Giants" },
    @{ Row = 41; Date = "01/04/2026"; Visitor = "Tennessee Titans"; Home = "Jacksonville Jaguars" },
    @{ Row = 42; Date = "01/04/2026"; Visitor = "Kansas City Chiefs"; Home = "Las Vegas Raiders" },
    @{ Row = 43; Date = "01/04/2026"; Visitor = "Arizona Cardinals"; Home = "Los Angeles Rams" },
    @{ Row = 44; Date = "01/04/2026"; Visitor = "Los Angeles Chargers"; Home = "Denver Broncos" },
    @{ Row = 45; Date = "01/04/2026"; Visitor = "Detroit Lions"; Home = "Chicago Bears" },
    @{ Row = 46; Date = "01/04/2026"; Visitor = "New York Jets"; Home = "Buffalo Bills" },
    @{ Row = 47; Date = "01/04/2026"; Visitor = "Miami Dolphins"; Home = "New England Patriots" },
    @{ Row = 48; Date = "01/04/2026"; Visitor = "Washington Commanders"; Home = "Philadelphia Eagles" },
    @{ Row = 49; Date = "01/04/2026"; Visitor = "Baltimore Ravens"; Home = "Pittsburgh Steelers" }
)

foreach ($r in $rows) {
    $dateCell = $ws.Range("B$($r.Row)")
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $r.Date
    $dateCell.Style = "Normal"
    $ws.Range("C$($r.Row)").Value = $r.Visitor
    $ws.Range("D$($r.Row)").Value = $r.Home
}
